$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new "name" cells first, in the order that yields the shared-string
# table layout produced by the original authoring session.
$ws.Range("B20").Value = "HiTechnicEOPD"
$ws.Range("B17").Value = "DexterPressureSensor250"
$ws.Range("B18").Value = "DexterPressureSensor500"
$ws.Range("F17").Value = "Pressure"
$ws.Range("F20").Value = "Distance"

# Row 17 - DPressure250 -> DexterPressureSensor250
$ws.Range("D17").Value = "Lawrie"
$ws.Range("E17").Value = "N"
$ws.Range("G17").Value = "SampleProvider"

# Row 18 - DPressure500 -> DexterPressureSensor500
$ws.Range("D18").Value = "Lawrie"
$ws.Range("E18").Value = "N"
$ws.Range("F18").Value = "Pressure"
$ws.Range("G18").Value = "SampleProvider"

# Row 20 - EOPD -> HiTechnicEOPD
$ws.Range("D20").Value = "Lawrie"
$ws.Range("E20").Value = "N"
$ws.Range("G20").Value = "SampleProvider"

# Move the active selection to G22, matching the author's last interaction
$ws.Range("G22").Select()
